# Rename the active sheet from "resultExcel" to "result"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "result"

# Clear the old header row (A1:B1) that held "Domain"/"Status"
$ws.Range("A1:B1").ClearContents()

# Write the new data row (row 6): domain check results
$ws.Range("A6").Value = "en.wikipedia.org"
$ws.Range("B6").Value = " UP "
$ws.Range("C6").Value = "200 = OK"
$ws.Range("D6").Value = "NO RESULT"

# Apply the same style (wrap text) used previously on A1/B1 to the new cells
$ws.Range("A6:D6").WrapText = $true

# Autofit columns to best-fit their content widths
$ws.Range("A1:D6").Columns.AutoFit()
